$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 74
$ws_ALC.Range("H74").Value = 3820
$ws_ALC.Range("I74").Value = 3297.5
$ws_ALC.Range("J74").Value = 4656
$ws_ALC.Range("K74").Value = 3297.5
$ws_ALC.Range("L74").Value = 4656
$ws_ALC.Range("M74").Value = -2361.5
$ws_ALC.Range("N74").Value = -6528

# ALC row 77
$ws_ALC.Range("H77").Value = 3820
$ws_ALC.Range("I77").Value = 3297.5
$ws_ALC.Range("J77").Value = 4656
$ws_ALC.Range("K77").Value = 16487.5
$ws_ALC.Range("L77").Value = 23280
$ws_ALC.Range("M77").Value = -11807.5
$ws_ALC.Range("N77").Value = -32640

# ALC row 112
$ws_ALC.Range("H112").Value = 1845.0476
$ws_ALC.Range("J112").Value = 1881.3684
$ws_ALC.Range("L112").Value = 5644.1052
$ws_ALC.Range("N112").Value = -7860.1052

# ALC row 137
$ws_ALC.Range("H137").Value = 5666.4443
$ws_ALC.Range("I137").Value = 10333
$ws_ALC.Range("J137").Value = 3333.1667
$ws_ALC.Range("K137").Value = 30999
$ws_ALC.Range("L137").Value = 9999.500100000001
$ws_ALC.Range("M137").Value = -28449
$ws_ALC.Range("N137").Value = -15099.5001

# ALC row 138
$ws_ALC.Range("H138").Value = 2003.5834
$ws_ALC.Range("I138").Value = 1963.0625
$ws_ALC.Range("J138").Value = 2023.8438
$ws_ALC.Range("K138").Value = 5889.1875
$ws_ALC.Range("L138").Value = 6071.5314
$ws_ALC.Range("M138").Value = -749.1875
$ws_ALC.Range("N138").Value = -16351.5314

# ARM row 61
$ws_ARM.Range("H61").Value = 2290.5715
$ws_ARM.Range("I61").Value = 1611.8462
$ws_ARM.Range("J61").Value = 11114
$ws_ARM.Range("K61").Value = 1611.8462
$ws_ARM.Range("L61").Value = 11114
$ws_ARM.Range("M61").Value = -1399.8462
$ws_ARM.Range("N61").Value = -11538

# ARM row 74
$ws_ARM.Range("H74").Value = 1941.6364
$ws_ARM.Range("I74").Value = 1967.8
$ws_ARM.Range("J74").Value = 1680
$ws_ARM.Range("K74").Value = 1967.8
$ws_ARM.Range("L74").Value = 1680
$ws_ARM.Range("M74").Value = -1093.8
$ws_ARM.Range("N74").Value = -3428

# ARM row 77
$ws_ARM.Range("H77").Value = 1941.6364
$ws_ARM.Range("I77").Value = 1967.8
$ws_ARM.Range("J77").Value = 1680
$ws_ARM.Range("K77").Value = 9839
$ws_ARM.Range("L77").Value = 8400
$ws_ARM.Range("M77").Value = -5471
$ws_ARM.Range("N77").Value = -17136

# ARM row 132
$ws_ARM.Range("H132").Value = 2402.34
$ws_ARM.Range("I132").Value = 1177.7407
$ws_ARM.Range("J132").Value = 3839.913
$ws_ARM.Range("K132").Value = 3533.2221
$ws_ARM.Range("L132").Value = 11519.739
$ws_ARM.Range("M132").Value = -1003.2221
$ws_ARM.Range("N132").Value = -16579.739

# ARM row 136
$ws_ARM.Range("H136").Value = 2290.5715
$ws_ARM.Range("I136").Value = 1611.8462
$ws_ARM.Range("J136").Value = 11114
$ws_ARM.Range("K136").Value = 4835.5386
$ws_ARM.Range("L136").Value = 33342
$ws_ARM.Range("M136").Value = -2285.5386
$ws_ARM.Range("N136").Value = -38442

# BSM row 107
$ws_BSM.Range("H107").Value = 1254.4667
$ws_BSM.Range("I107").Value = 890.84
$ws_BSM.Range("J107").Value = 3072.6
$ws_BSM.Range("K107").Value = 890.84
$ws_BSM.Range("L107").Value = 3072.6
$ws_BSM.Range("M107").Value = 1029.16
$ws_BSM.Range("N107").Value = -6912.6

# BSM row 134
$ws_BSM.Range("H134").Value = 2477.8386
$ws_BSM.Range("I134").Value = 1174.25
$ws_BSM.Range("J134").Value = 4848
$ws_BSM.Range("K134").Value = 3522.75
$ws_BSM.Range("L134").Value = 14544
$ws_BSM.Range("M134").Value = -987.75
$ws_BSM.Range("N134").Value = -19614

# CRP row 31
$ws_CRP.Range("H31").Value = 2246.8696
$ws_CRP.Range("I31").Value = 1954.7931
$ws_CRP.Range("J31").Value = 2745.1177
$ws_CRP.Range("K31").Value = 1954.7931
$ws_CRP.Range("L31").Value = 2745.1177
$ws_CRP.Range("M31").Value = -1659.7931
$ws_CRP.Range("N31").Value = -3335.1177

# CRP row 34
$ws_CRP.Range("H34").Value = 2246.8696
$ws_CRP.Range("I34").Value = 1954.7931
$ws_CRP.Range("J34").Value = 2745.1177
$ws_CRP.Range("K34").Value = 1954.7931
$ws_CRP.Range("L34").Value = 2745.1177
$ws_CRP.Range("M34").Value = -1752.7931
$ws_CRP.Range("N34").Value = -3149.1177

# CRP row 58
$ws_CRP.Range("H58").Value = 2168.5676
$ws_CRP.Range("I58").Value = 1718.1666
$ws_CRP.Range("J58").Value = 2384.76
$ws_CRP.Range("K58").Value = 1718.1666
$ws_CRP.Range("L58").Value = 2384.76
$ws_CRP.Range("M58").Value = -1515.1666
$ws_CRP.Range("N58").Value = -2790.76

# CRP row 132
$ws_CRP.Range("H132").Value = 2454.7856
$ws_CRP.Range("I132").Value = 1351.75
$ws_CRP.Range("K132").Value = 4055.25
$ws_CRP.Range("M132").Value = -1525.25

# CRP row 134
$ws_CRP.Range("H134").Value = 2888.25
$ws_CRP.Range("I134").Value = 3316.3333
$ws_CRP.Range("J134").Value = 2246.125
$ws_CRP.Range("K134").Value = 9948.999899999999
$ws_CRP.Range("L134").Value = 6738.375
$ws_CRP.Range("M134").Value = -7413.999899999999
$ws_CRP.Range("N134").Value = -11808.375

# CRP row 136
$ws_CRP.Range("H136").Value = 2168.5676
$ws_CRP.Range("I136").Value = 1718.1666
$ws_CRP.Range("J136").Value = 2384.76
$ws_CRP.Range("K136").Value = 5154.4998
$ws_CRP.Range("L136").Value = 7154.280000000001
$ws_CRP.Range("M136").Value = -2604.4998
$ws_CRP.Range("N136").Value = -12254.28

# CUL row 107
$ws_CUL.Range("H107").Value = 1338.25
$ws_CUL.Range("I107").Value = 2927.25
$ws_CUL.Range("J107").Value = 543.75
$ws_CUL.Range("K107").Value = 8781.75
$ws_CUL.Range("L107").Value = 1631.25
$ws_CUL.Range("M107").Value = -6861.75
$ws_CUL.Range("N107").Value = -5471.25

# CUL row 132
$ws_CUL.Range("H132").Value = 1371.9546
$ws_CUL.Range("I132").Value = 1345.4
$ws_CUL.Range("J132").Value = 1394.0834
$ws_CUL.Range("K132").Value = 12108.6
$ws_CUL.Range("L132").Value = 12546.7506
$ws_CUL.Range("M132").Value = -9578.6
$ws_CUL.Range("N132").Value = -17606.7506

# GSM row 132
$ws_GSM.Range("H132").Value = 3198.4583
$ws_GSM.Range("I132").Value = 2294.3125
$ws_GSM.Range("J132").Value = 5006.75
$ws_GSM.Range("K132").Value = 6882.9375
$ws_GSM.Range("L132").Value = 15020.25
$ws_GSM.Range("M132").Value = -4352.9375
$ws_GSM.Range("N132").Value = -20080.25

# LTW row 132
$ws_LTW.Range("H132").Value = 8353.82
$ws_LTW.Range("I132").Value = 11382.091
$ws_LTW.Range("K132").Value = 34146.273
$ws_LTW.Range("M132").Value = -31616.273

# LTW row 136
$ws_LTW.Range("H136").Value = 9806824
$ws_LTW.Range("I136").Value = 2877.111
$ws_LTW.Range("J136").Value = 20836264
$ws_LTW.Range("K136").Value = 8631.332999999999
$ws_LTW.Range("L136").Value = 62508792
$ws_LTW.Range("M136").Value = -6081.332999999999
$ws_LTW.Range("N136").Value = -62513892

# WVR row 81
$ws_WVR.Range("H81").Value = 1180.25
$ws_WVR.Range("I81").Value = 1005
$ws_WVR.Range("J81").Value = 1355.5
$ws_WVR.Range("K81").Value = 2010
$ws_WVR.Range("L81").Value = 2711
$ws_WVR.Range("M81").Value = -949
$ws_WVR.Range("N81").Value = -4833

# WVR row 84
$ws_WVR.Range("H84").Value = 1180.25
$ws_WVR.Range("I84").Value = 1005
$ws_WVR.Range("J84").Value = 1355.5
$ws_WVR.Range("K84").Value = 10050
$ws_WVR.Range("L84").Value = 13555
$ws_WVR.Range("M84").Value = -4746
$ws_WVR.Range("N84").Value = -24163

# WVR row 132
$ws_WVR.Range("H132").Value = 2117.6
$ws_WVR.Range("I132").Value = 1650.2174
$ws_WVR.Range("J132").Value = 3653.2856
$ws_WVR.Range("K132").Value = 4950.6522
$ws_WVR.Range("L132").Value = 10959.8568
$ws_WVR.Range("M132").Value = -2420.6522
$ws_WVR.Range("N132").Value = -16019.8568

# WVR row 136
$ws_WVR.Range("H136").Value = 3021.6177
$ws_WVR.Range("I136").Value = 810.7143
$ws_WVR.Range("J136").Value = 6593.077
$ws_WVR.Range("K136").Value = 2432.1429
$ws_WVR.Range("L136").Value = 19779.231
$ws_WVR.Range("M136").Value = 117.8571000000002
$ws_WVR.Range("N136").Value = -24879.231
